$d = $word.ActiveDocument

# --- Edit 1 -------------------------------------------------------------
# "... text is removed (e.g., of citations."
#   -> "... text is removed in citations."
$r1 = $d.Content
$found1 = $r1.Find.Execute(" (e.g., of citation", $true, $false, $false, $false, $false, $true, 1, $false, " in citation", 2)
Write-Output "Edit1 found: $found1"
if (-not $found1) {
    throw "Edit1: target text ' (e.g., of citation' was not found"
}

# --- Edit 2 -------------------------------------------------------------
# "... periods (.), and commas ..."  ->  "... periods (,), and commas ..."
$r2 = $d.Content
$found2 = $r2.Find.Execute("periods (.), and commas", $true, $false, $false, $false, $false, $true, 1, $false, "periods (,), and commas", 2)
Write-Output "Edit2 found: $found2"
if (-not $found2) {
    throw "Edit2: target text 'periods (.), and commas' was not found"
}
